$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.937.85"
$ws.Range("E2").Value = "  +0.96%  "

$ws.Range("D3").Value = "1.882.10"
$ws.Range("E3").Value = "  +0.16%  "

$ws.Range("E4").Value = "  +1.77%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.97%  "

$ws.Range("E6").Value = "  +1.65%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4647"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.78%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3895"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.94%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.88"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.86%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07909"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.73%  "

$ws.Range("E11").Value = "  -2.00%  "

$ws.Range("E12").Value = "  -1.32%  "

$ws.Range("D13").Value = "1.898.48"
$ws.Range("E13").Value = "  +1.89%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.907"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.00%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.060"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.64%  "

$ws.Range("E16").Value = "  +1.76%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.06735"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.73%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "86.54"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.59%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001038"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.45%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.89%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.018"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.56%  "

$ws.Range("D22").Value = "27.963.95"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.444"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.27%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.38%  "

$ws.Range("E25").Value = "  +2.49%  "

$ws.Range("D26").Value = "2.112.47"
$ws.Range("E26").Value = "  +1.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.57"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.11%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.82"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.19%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.047"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.22%  "

$ws.Range("E30").Value = "  -3.50%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "120.81"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.31%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09447"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.34%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9527"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.32%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.664"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.08%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.290"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.52%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.341"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.58%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06070"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.83%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02225"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.203"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.31%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.072"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.84%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5897"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.01%  "

$ws.Range("E42").Value = "  -1.67%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.58%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.270"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.14%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5606"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.02%  "

$ws.Range("E46").Value = "  -1.89%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.378"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.85%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.899"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.81%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06893"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.03%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "113.03"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.65%  "

$ws.Range("E51").Value = "  -1.29%  "
